# Insert a new column A (shifts all existing columns A:AA -> B:AB)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Fix new column A's width (target rendered width "23")
$ws.Columns("A:A").ColumnWidth = 22.1666666666667

# Give column A (rows 1-17) the same formatting (style) as column B
$ws.Range("B1:B17").Copy()
$ws.Range("A1:A17").PasteSpecial(-4122)

# Header row text (uppercased + new INDEX column)
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("B1").Value = "CATEGORY"
$ws.Range("C1").Value = "REGION"
$ws.Range("D1").Value = "DIVISION"
$ws.Range("E1").Value = "SCHOOL ID"
$ws.Range("F1").Value = "SCHOOL NAME"
$ws.Range("G1").Value = "MUNICIPALITY"
$ws.Range("H1").Value = "LD"
$ws.Range("I1").Value = "TOTAL NO. OF SITES"
$ws.Range("J1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("K1").Value = "SCOPE OF WORK"
$ws.Range("L1").Value = "PROJECT ALLOCATION"
$ws.Range("M1").Value = "BATCH"
$ws.Range("N1").Value = "CONTRACT AMOUNT"
$ws.Range("O1").Value = "STATUS"
$ws.Range("P1").Value = "PERCENTAGE OF COMPLETION"
$ws.Range("Q1").Value = " TARGET COMPLETION DATE "
$ws.Range("R1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("S1").Value = "PROJECT ID"
$ws.Range("T1").Value = "CONTRACT ID"
$ws.Range("U1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("V1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("W1").Value = "BID OPENING"
$ws.Range("X1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("Y1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("Z1").Value = "NAME OF CONTRACTOR"
$ws.Range("AA1").Value = "OTHER REMARKS"
$ws.Range("AB1").Value = "Status as of July 11, 2025"

# New INDEX values for data rows (2-17): sequential 29..44
for ($i = 2; $i -le 17; $i++) {
    $ws.Cells.Item($i, 1).Value = 27 + $i
}
